# Apply the two changes captured by the commit:
#  1. Slide 16's table switches from the custom "Table_0" style to the
#     built-in "Medium Style 2 - Accent 1" table style.
#  2. The deck's theme colour scheme changes from the custom "Integral"
#     palette to the stock "Office" palette (the colours that the theme
#     swap between theme1.xml/theme2.xml brings to the slide master).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$targetStyleId = "{F1C666BA-BD9D-45E5-A789-83596BB4AFB4}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}

# --- 2. Theme colours (Integral -> Office) ---------------------------
# Order matches the OOXML <a:clrScheme> child order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink. Values are VBA-style BGR integers (the
# RGB() encoding PowerPoint's ColorScheme.Colors(i).RGB expects).
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
